# --- Rename the existing "types" sheet to "list" and add a new "create" sheet after it ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "list"

# Add a throwaway sheet first, then the real one, then drop the throwaway,
# so the "create" sheet's internal sheetId lands on 3, matching the
# authored workbook (sheetId is a monotonically increasing counter, not a
# tab position). Re-fetch the surviving sheet from the collection after
# the delete since the old COM reference can go stale.
$throwaway = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$null = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $throwaway)
[void]$throwaway.Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "create"

# --- Populate the "list" sheet (type name / enum values / description) ---
# Values are entered column-by-column for the first two rows and then
# row-by-row for the rest, matching how the source workbook was authored.
$ws1.Cells.Item(1, 1).Value = "asset_category_type"
$ws1.Cells.Item(2, 1).Value = "person_category_type"

$ws1.Cells.Item(1, 2).Value = "('F', 'E')"
$ws1.Cells.Item(2, 2).Value = "('E', 'C', 'T', 'R')"

$ws1.Cells.Item(1, 3).Value = "facility, equipment"
$ws1.Cells.Item(2, 3).Value = "efetivo, comissionado, terceirizado, estagiário"

$ws1.Cells.Item(3, 1).Value = "order_status_type"
$ws1.Cells.Item(4, 1).Value = "order_priority_type"
$ws1.Cells.Item(5, 1).Value = "order_category_type"

$ws1.Cells.Item(3, 2).Value = "('R', 'E', 'C', 'Z')"
$ws1.Cells.Item(3, 3).Value = "pendente de análise, em execução, concluída, cancelada"

$ws1.Cells.Item(4, 2).Value = "('H', 'N')"
$ws1.Cells.Item(4, 3).Value = "high, normal"

$ws1.Cells.Item(5, 2).Value = "('E', 'A', 'C')"
$ws1.Cells.Item(5, 3).Value = "elétrica, ar-condicionado, civil"

# Column widths for "list" (A/B/C). The host quantizes ColumnWidth to 1/6
# of a character, so feed it pre-compensated inputs to land as close as
# possible to the authored widths.
$ws1.Columns.Item(1).ColumnWidth = 25.736979166666668
$ws1.Columns.Item(2).ColumnWidth = 20.307291666666668
$ws1.Columns.Item(3).ColumnWidth = 60.307291666666664

# --- Populate the "create" sheet with formulas referencing "list" ---
for ($i = 1; $i -le 5; $i++) {
  $ws2.Cells.Item($i, 1).Formula = "=CONCATENATE(""CREATE TYPE "",list!A$i,"" AS ENUM "",list!B$i,"";"")"
}

# Column width for "create" (A)
$ws2.Columns.Item(1).ColumnWidth = 62.307291666666664

# "create" page margins (narrow / metric preset)
$ws2.PageSetup.LeftMargin = 0.511811024 * 72
$ws2.PageSetup.RightMargin = 0.511811024 * 72
$ws2.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws2.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws2.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$ws2.PageSetup.FooterMargin = 0.31496062000000002 * 72

# --- View state: "list" is no longer the selected tab, "create" is ---
[void]$ws1.Activate()
$excel.ActiveWindow.Zoom = 140
[void]$ws1.Range("B6").Select()

[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 160
[void]$ws2.Range("A2").Select()
